$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format so numeric-looking strings
# (e.g. "0.9988", "246.62") are stored as text, matching the source data
# which uses inline strings for every Price/Volume cell.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Rows 2-33: coin identity unchanged, update Price / Volume(1h) ---
$ws.Range("D2").Value = "26.604.51"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "1.730.73"

$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "246.62"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "0.4817"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.2680"
$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("D10").Value = "1.732.29"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").Value = "0.07150"
$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D12").Value = "15.65"
$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("D13").Value = "0.6129"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").Value = "4.547"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").Value = "77.41"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "0.9993"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "26.601.11"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "0.9989"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").Value = "0.000006968"
$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("D20").Value = "11.58"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("D21").Value = "1.951.70"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").Value = "4.522"
$ws.Range("E22").Value = "  -2.89%  "

$ws.Range("D23").Value = "8.816"
$ws.Range("E23").Value = "  -0.56%  "

$ws.Range("D24").Value = "5.249"
$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("D25").Value = "137.26"
$ws.Range("E25").Value = "  +0.80%  "

$ws.Range("D26").Value = "15.42"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").Value = "1.785"
$ws.Range("E27").Value = "  -2.22%  "

$ws.Range("D28").Value = "1.411"
$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").Value = "108.60"
$ws.Range("E29").Value = "  +0.74%  "

$ws.Range("D30").Value = "3.981"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("D31").Value = "0.08023"
$ws.Range("E31").Value = "  +1.56%  "

$ws.Range("D32").Value = "3.693"
$ws.Range("E32").Value = "  -2.13%  "

$ws.Range("D33").Value = "0.04543"
$ws.Range("E33").Value = "  -0.86%  "

# --- Rows 34-51: new "Frax" entry inserted, subsequent coins shift down one rank,
#     the previous last entry ("Aave") drops off the bottom of the table ---
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "0.9988"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.614"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.6350"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.055"
$ws.Range("E38").Value = "  +3.81%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "0.8984"
$ws.Range("E39").Value = "  -6.33%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.373"
$ws.Range("E40").Value = "  -3.99%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "103.19"
$ws.Range("E42").Value = "  -10.03%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.01504"
$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.474"
$ws.Range("E44").Value = "  -4.27%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "7.200"
$ws.Range("E45").Value = "  +6.61%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "0.3910"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1188"
$ws.Range("E47").Value = "  -1.52%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05384"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.897"
$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "30.73"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.256"
$ws.Range("E51").Value = "  +0.56%  "
